$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value would otherwise be auto-parsed as a number by Excel;
# the source data keeps the Price column as text, so force text format before writing.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '65.739.77'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '2.674.76'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '600.19'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").Value = '155.93'
$ws.Range("E6").Value = '  -1.67%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.603'
$ws.Range("E8").Value = '  +2.04%  '
$ws.Range("E9").Value = '  -1.19%  '
$ws.Range("E10").Value = '  +1.48%  '
$ws.Range("D11").Value = '0.397'
$ws.Range("E11").Value = '  -2.26%  '
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").Value = '29.37'
$ws.Range("E13").Value = '  -1.54%  '
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("D15").Value = '3.153.73'
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '65.508.76'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D17").Value = '2.679.83'
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").Value = '12.52'
$ws.Range("E18").Value = '  -1.95%  '
$ws.Range("E19").Value = '  -1.84%  '
$ws.Range("E20").Value = '  +1.67%  '
$ws.Range("D21").Value = '350.58'
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").Value = '70.15'
$ws.Range("E23").Value = '  +1.68%  '
$ws.Range("D24").Value = '9.77'
$ws.Range("E24").Value = '  +1.96%  '
$ws.Range("E25").Value = '  +2.64%  '
$ws.Range("E26").Value = '  -3.92%  '
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("D28").Value = '0.169'
$ws.Range("E28").Value = '  +1.78%  '
$ws.Range("D29").Value = '8.10'
$ws.Range("E29").Value = '  -1.77%  '
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("D31").Value = '538.20'
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("E32").Value = '  -2.74%  '
$ws.Range("D33").Value = '1.76'
$ws.Range("E33").Value = '  -5.40%  '
$ws.Range("E34").Value = '  +2.51%  '
$ws.Range("E35").Value = '  -4.00%  '
$ws.Range("D36").Value = '0.424'
$ws.Range("E36").Value = '  -2.57%  '
$ws.Range("D37").Value = '20.37'
$ws.Range("E37").Value = '  -1.56%  '
$ws.Range("D38").Value = '159.53'
$ws.Range("E38").Value = '  -2.08%  '
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("E40").Value = '  -3.66%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D42").Value = '42.45'
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").Value = '166.18'
$ws.Range("E43").Value = '  -0.42%  '
$ws.Range("D44").Value = '4.09'
$ws.Range("E44").Value = '  -2.38%  '
$ws.Range("D45").Value = '0.0611'
$ws.Range("E45").Value = '  -0.53%  '
$ws.Range("D46").Value = '22.96'
$ws.Range("E46").Value = '  -1.06%  '
$ws.Range("D47").Value = '2.24'
$ws.Range("E47").Value = '  -5.52%  '
$ws.Range("D48").Value = '0.647'
$ws.Range("E48").Value = '  -2.39%  '
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("E50").Value = '  +0.56%  '
$ws.Range("D51").Value = '19.98'
$ws.Range("E51").Value = '  +0.57%  '
